$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Components")

# Insert new column M for pin_names (shifts hfe/speed to N/O)
$ws.Columns.Item(13).Insert()

# Header row
$ws.Cells.Item(1,1).Value = 'id'
$ws.Cells.Item(1,2).Value = 'category'
$ws.Cells.Item(1,3).Value = 'polarity'
$ws.Cells.Item(1,4).Value = 'package_id'
$ws.Cells.Item(1,5).Value = 'pinout_code'
$ws.Cells.Item(1,6).Value = 'v_max'
$ws.Cells.Item(1,7).Value = 'i_max'
$ws.Cells.Item(1,8).Value = 'power_max'
$ws.Cells.Item(1,9).Value = 'v_trig'
$ws.Cells.Item(1,10).Value = 'r_ds'
$ws.Cells.Item(1,11).Value = 'test_script_id'
$ws.Cells.Item(1,12).Value = 'description'
$ws.Cells.Item(1,13).Value = 'pin_names'
$ws.Cells.Item(1,14).Value = 'hfe'
$ws.Cells.Item(1,15).Value = 'speed'

# Data rows
# Row 2: IRFZ44N
$ws.Cells.Item(2,1).Value = 'IRFZ44N'
$ws.Cells.Item(2,2).Value = 'MOSFET'
$ws.Cells.Item(2,3).Value = 'N-Channel'
$ws.Cells.Item(2,4).Value = 'TO-220'
$ws.Cells.Item(2,5).Value = 'GDS'
$ws.Cells.Item(2,6).Value = 55
$ws.Cells.Item(2,7).Value = 49
$ws.Cells.Item(2,8).Value = 94
$ws.Cells.Item(2,9).Value = 4
$ws.Cells.Item(2,10).Value = 0.017
$ws.Cells.Item(2,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(2,12).Value = 'Standart N-Kanal Güç MOSFET''i. Motor sürücü ve güç kaynaklarında yaygın.'
$ws.Cells.Item(2,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(2,14).Value = ''
$ws.Cells.Item(2,15).Value = ''

# Row 3: IRF3205
$ws.Cells.Item(3,1).Value = 'IRF3205'
$ws.Cells.Item(3,2).Value = 'MOSFET'
$ws.Cells.Item(3,3).Value = 'N-Channel'
$ws.Cells.Item(3,4).Value = 'TO-220'
$ws.Cells.Item(3,5).Value = 'GDS'
$ws.Cells.Item(3,6).Value = 55
$ws.Cells.Item(3,7).Value = 110
$ws.Cells.Item(3,8).Value = 200
$ws.Cells.Item(3,9).Value = 4
$ws.Cells.Item(3,10).Value = 0.008
$ws.Cells.Item(3,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(3,12).Value = 'Çok düşük iç dirençli, yüksek akım MOSFET''i. İnverterler için ideal.'
$ws.Cells.Item(3,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(3,14).Value = ''
$ws.Cells.Item(3,15).Value = ''

# Row 4: IRF540N
$ws.Cells.Item(4,1).Value = 'IRF540N'
$ws.Cells.Item(4,2).Value = 'MOSFET'
$ws.Cells.Item(4,3).Value = 'N-Channel'
$ws.Cells.Item(4,4).Value = 'TO-220'
$ws.Cells.Item(4,5).Value = 'GDS'
$ws.Cells.Item(4,6).Value = 100
$ws.Cells.Item(4,7).Value = 33
$ws.Cells.Item(4,8).Value = 130
$ws.Cells.Item(4,9).Value = 4
$ws.Cells.Item(4,10).Value = 0.044
$ws.Cells.Item(4,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(4,12).Value = '100V Dayanımlı genel amaçlı güç MOSFET''i.'
$ws.Cells.Item(4,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(4,14).Value = ''
$ws.Cells.Item(4,15).Value = ''

# Row 5: IRF640
$ws.Cells.Item(5,1).Value = 'IRF640'
$ws.Cells.Item(5,2).Value = 'MOSFET'
$ws.Cells.Item(5,3).Value = 'N-Channel'
$ws.Cells.Item(5,4).Value = 'TO-220'
$ws.Cells.Item(5,5).Value = 'GDS'
$ws.Cells.Item(5,6).Value = 200
$ws.Cells.Item(5,7).Value = 18
$ws.Cells.Item(5,8).Value = 150
$ws.Cells.Item(5,9).Value = 4
$ws.Cells.Item(5,10).Value = 0.15
$ws.Cells.Item(5,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(5,12).Value = '200V Yüksek voltaj anahtarlama elemanı.'
$ws.Cells.Item(5,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(5,14).Value = ''
$ws.Cells.Item(5,15).Value = ''

# Row 6: IRF740
$ws.Cells.Item(6,1).Value = 'IRF740'
$ws.Cells.Item(6,2).Value = 'MOSFET'
$ws.Cells.Item(6,3).Value = 'N-Channel'
$ws.Cells.Item(6,4).Value = 'TO-220'
$ws.Cells.Item(6,5).Value = 'GDS'
$ws.Cells.Item(6,6).Value = 400
$ws.Cells.Item(6,7).Value = 10
$ws.Cells.Item(6,8).Value = 125
$ws.Cells.Item(6,9).Value = 4
$ws.Cells.Item(6,10).Value = 0.55
$ws.Cells.Item(6,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(6,12).Value = '400V SMPS ve sürücü devreleri için.'
$ws.Cells.Item(6,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(6,14).Value = ''
$ws.Cells.Item(6,15).Value = ''

# Row 7: IRF840
$ws.Cells.Item(7,1).Value = 'IRF840'
$ws.Cells.Item(7,2).Value = 'MOSFET'
$ws.Cells.Item(7,3).Value = 'N-Channel'
$ws.Cells.Item(7,4).Value = 'TO-220'
$ws.Cells.Item(7,5).Value = 'GDS'
$ws.Cells.Item(7,6).Value = 500
$ws.Cells.Item(7,7).Value = 8
$ws.Cells.Item(7,8).Value = 125
$ws.Cells.Item(7,9).Value = 4
$ws.Cells.Item(7,10).Value = 0.85
$ws.Cells.Item(7,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(7,12).Value = '500V Yüksek voltaj uygulamaları.'
$ws.Cells.Item(7,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(7,14).Value = ''
$ws.Cells.Item(7,15).Value = ''

# Row 8: IRF9540
$ws.Cells.Item(8,1).Value = 'IRF9540'
$ws.Cells.Item(8,2).Value = 'MOSFET'
$ws.Cells.Item(8,3).Value = 'P-Channel'
$ws.Cells.Item(8,4).Value = 'TO-220'
$ws.Cells.Item(8,5).Value = 'GDS'
$ws.Cells.Item(8,6).Value = -100
$ws.Cells.Item(8,7).Value = -23
$ws.Cells.Item(8,8).Value = 140
$ws.Cells.Item(8,9).Value = -4
$ws.Cells.Item(8,10).Value = 0.117
$ws.Cells.Item(8,11).Value = 'TEST_MOSFET_P'
$ws.Cells.Item(8,12).Value = 'P-Kanal Güç MOSFET''i (High-Side Anahtarlama).'
$ws.Cells.Item(8,13).Value = 'GATE,DRAIN,SOURCE'
$ws.Cells.Item(8,14).Value = ''
$ws.Cells.Item(8,15).Value = ''

# Row 9: 2N7000
$ws.Cells.Item(9,1).Value = '2N7000'
$ws.Cells.Item(9,2).Value = 'MOSFET'
$ws.Cells.Item(9,3).Value = 'N-Channel'
$ws.Cells.Item(9,4).Value = 'TO-92'
$ws.Cells.Item(9,5).Value = 'SGD'
$ws.Cells.Item(9,6).Value = 60
$ws.Cells.Item(9,7).Value = 0.2
$ws.Cells.Item(9,8).Value = 0.4
$ws.Cells.Item(9,9).Value = 2.1
$ws.Cells.Item(9,10).Value = 5
$ws.Cells.Item(9,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(9,12).Value = 'Küçük sinyal, lojik seviye MOSFET. Arduino ile sürülebilir.'
$ws.Cells.Item(9,13).Value = 'SOURCE,GATE,DRAIN'
$ws.Cells.Item(9,14).Value = ''
$ws.Cells.Item(9,15).Value = ''

# Row 10: BS170
$ws.Cells.Item(10,1).Value = 'BS170'
$ws.Cells.Item(10,2).Value = 'MOSFET'
$ws.Cells.Item(10,3).Value = 'N-Channel'
$ws.Cells.Item(10,4).Value = 'TO-92'
$ws.Cells.Item(10,5).Value = 'DGS'
$ws.Cells.Item(10,6).Value = 60
$ws.Cells.Item(10,7).Value = 0.5
$ws.Cells.Item(10,8).Value = 0.8
$ws.Cells.Item(10,9).Value = 2.1
$ws.Cells.Item(10,10).Value = 5
$ws.Cells.Item(10,11).Value = 'TEST_MOSFET_N'
$ws.Cells.Item(10,12).Value = '2N7000 benzeri ama bacak dizilimi farklı (DGS).'
$ws.Cells.Item(10,13).Value = 'DRAIN,GATE,SOURCE'
$ws.Cells.Item(10,14).Value = ''
$ws.Cells.Item(10,15).Value = ''

# Row 11: BC547
$ws.Cells.Item(11,1).Value = 'BC547'
$ws.Cells.Item(11,2).Value = 'BJT'
$ws.Cells.Item(11,3).Value = 'NPN'
$ws.Cells.Item(11,4).Value = 'TO-92'
$ws.Cells.Item(11,5).Value = 'CBE'
$ws.Cells.Item(11,6).Value = 45
$ws.Cells.Item(11,7).Value = 0.1
$ws.Cells.Item(11,8).Value = 0.5
$ws.Cells.Item(11,9).Value = ''
$ws.Cells.Item(11,10).Value = ''
$ws.Cells.Item(11,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(11,12).Value = 'Genel amaçlı, düşük gürültülü NPN.'
$ws.Cells.Item(11,13).Value = 'COLL,BASE,EMIT'
$ws.Cells.Item(11,14).Value = 110
$ws.Cells.Item(11,15).Value = ''

# Row 12: BC557
$ws.Cells.Item(12,1).Value = 'BC557'
$ws.Cells.Item(12,2).Value = 'BJT'
$ws.Cells.Item(12,3).Value = 'PNP'
$ws.Cells.Item(12,4).Value = 'TO-92'
$ws.Cells.Item(12,5).Value = 'CBE'
$ws.Cells.Item(12,6).Value = -45
$ws.Cells.Item(12,7).Value = -0.1
$ws.Cells.Item(12,8).Value = 0.5
$ws.Cells.Item(12,9).Value = ''
$ws.Cells.Item(12,10).Value = ''
$ws.Cells.Item(12,11).Value = 'TEST_BJT_PNP'
$ws.Cells.Item(12,12).Value = 'Genel amaçlı PNP transistör.'
$ws.Cells.Item(12,13).Value = 'COLL,BASE,EMIT'
$ws.Cells.Item(12,14).Value = 110
$ws.Cells.Item(12,15).Value = ''

# Row 13: 2N2222
$ws.Cells.Item(13,1).Value = '2N2222'
$ws.Cells.Item(13,2).Value = 'BJT'
$ws.Cells.Item(13,3).Value = 'NPN'
$ws.Cells.Item(13,4).Value = 'TO-92'
$ws.Cells.Item(13,5).Value = 'EBC'
$ws.Cells.Item(13,6).Value = 40
$ws.Cells.Item(13,7).Value = 0.8
$ws.Cells.Item(13,8).Value = 0.5
$ws.Cells.Item(13,9).Value = ''
$ws.Cells.Item(13,10).Value = ''
$ws.Cells.Item(13,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(13,12).Value = 'Yüksek hızlı anahtarlama ve yükseltme.'
$ws.Cells.Item(13,13).Value = 'EMIT,BASE,COLL'
$ws.Cells.Item(13,14).Value = 100
$ws.Cells.Item(13,15).Value = ''

# Row 14: 2N3904
$ws.Cells.Item(14,1).Value = '2N3904'
$ws.Cells.Item(14,2).Value = 'BJT'
$ws.Cells.Item(14,3).Value = 'NPN'
$ws.Cells.Item(14,4).Value = 'TO-92'
$ws.Cells.Item(14,5).Value = 'EBC'
$ws.Cells.Item(14,6).Value = 40
$ws.Cells.Item(14,7).Value = 0.2
$ws.Cells.Item(14,8).Value = 0.6
$ws.Cells.Item(14,9).Value = ''
$ws.Cells.Item(14,10).Value = ''
$ws.Cells.Item(14,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(14,12).Value = 'Genel amaçlı NPN.'
$ws.Cells.Item(14,13).Value = 'EMIT,BASE,COLL'
$ws.Cells.Item(14,14).Value = 100
$ws.Cells.Item(14,15).Value = ''

# Row 15: BD139
$ws.Cells.Item(15,1).Value = 'BD139'
$ws.Cells.Item(15,2).Value = 'BJT'
$ws.Cells.Item(15,3).Value = 'NPN'
$ws.Cells.Item(15,4).Value = 'TO-126'
$ws.Cells.Item(15,5).Value = 'ECB'
$ws.Cells.Item(15,6).Value = 80
$ws.Cells.Item(15,7).Value = 1.5
$ws.Cells.Item(15,8).Value = 12
$ws.Cells.Item(15,9).Value = ''
$ws.Cells.Item(15,10).Value = ''
$ws.Cells.Item(15,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(15,12).Value = 'Orta güç NPN. Ses sürücü devrelerinde sıkça kullanılır.'
$ws.Cells.Item(15,13).Value = 'EMIT,COLL,BASE'
$ws.Cells.Item(15,14).Value = 63
$ws.Cells.Item(15,15).Value = ''

# Row 16: BD140
$ws.Cells.Item(16,1).Value = 'BD140'
$ws.Cells.Item(16,2).Value = 'BJT'
$ws.Cells.Item(16,3).Value = 'PNP'
$ws.Cells.Item(16,4).Value = 'TO-126'
$ws.Cells.Item(16,5).Value = 'ECB'
$ws.Cells.Item(16,6).Value = -80
$ws.Cells.Item(16,7).Value = -1.5
$ws.Cells.Item(16,8).Value = 12
$ws.Cells.Item(16,9).Value = ''
$ws.Cells.Item(16,10).Value = ''
$ws.Cells.Item(16,11).Value = 'TEST_BJT_PNP'
$ws.Cells.Item(16,12).Value = 'BD139''un PNP eşleniği.'
$ws.Cells.Item(16,13).Value = 'EMIT,COLL,BASE'
$ws.Cells.Item(16,14).Value = 63
$ws.Cells.Item(16,15).Value = ''

# Row 17: TIP31C
$ws.Cells.Item(17,1).Value = 'TIP31C'
$ws.Cells.Item(17,2).Value = 'BJT'
$ws.Cells.Item(17,3).Value = 'NPN'
$ws.Cells.Item(17,4).Value = 'TO-220'
$ws.Cells.Item(17,5).Value = 'BCE'
$ws.Cells.Item(17,6).Value = 100
$ws.Cells.Item(17,7).Value = 3
$ws.Cells.Item(17,8).Value = 40
$ws.Cells.Item(17,9).Value = ''
$ws.Cells.Item(17,10).Value = ''
$ws.Cells.Item(17,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(17,12).Value = 'Güç NPN Transistörü.'
$ws.Cells.Item(17,13).Value = 'BASE,COLL,EMIT'
$ws.Cells.Item(17,14).Value = 25
$ws.Cells.Item(17,15).Value = ''

# Row 18: TIP32C
$ws.Cells.Item(18,1).Value = 'TIP32C'
$ws.Cells.Item(18,2).Value = 'BJT'
$ws.Cells.Item(18,3).Value = 'PNP'
$ws.Cells.Item(18,4).Value = 'TO-220'
$ws.Cells.Item(18,5).Value = 'BCE'
$ws.Cells.Item(18,6).Value = -100
$ws.Cells.Item(18,7).Value = -3
$ws.Cells.Item(18,8).Value = 40
$ws.Cells.Item(18,9).Value = ''
$ws.Cells.Item(18,10).Value = ''
$ws.Cells.Item(18,11).Value = 'TEST_BJT_PNP'
$ws.Cells.Item(18,12).Value = 'Güç PNP Transistörü.'
$ws.Cells.Item(18,13).Value = 'BASE,COLL,EMIT'
$ws.Cells.Item(18,14).Value = 25
$ws.Cells.Item(18,15).Value = ''

# Row 19: TIP120
$ws.Cells.Item(19,1).Value = 'TIP120'
$ws.Cells.Item(19,2).Value = 'BJT'
$ws.Cells.Item(19,3).Value = 'NPN Darlington'
$ws.Cells.Item(19,4).Value = 'TO-220'
$ws.Cells.Item(19,5).Value = 'BCE'
$ws.Cells.Item(19,6).Value = 60
$ws.Cells.Item(19,7).Value = 5
$ws.Cells.Item(19,8).Value = 65
$ws.Cells.Item(19,9).Value = ''
$ws.Cells.Item(19,10).Value = ''
$ws.Cells.Item(19,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(19,12).Value = 'Darlington NPN. Çok yüksek kazançlı.'
$ws.Cells.Item(19,13).Value = 'BASE,COLL,EMIT'
$ws.Cells.Item(19,14).Value = 1000
$ws.Cells.Item(19,15).Value = ''

# Row 20: 2N3055
$ws.Cells.Item(20,1).Value = '2N3055'
$ws.Cells.Item(20,2).Value = 'BJT'
$ws.Cells.Item(20,3).Value = 'NPN Power'
$ws.Cells.Item(20,4).Value = 'TO-3'
$ws.Cells.Item(20,5).Value = 'BCE'
$ws.Cells.Item(20,6).Value = 60
$ws.Cells.Item(20,7).Value = 15
$ws.Cells.Item(20,8).Value = 115
$ws.Cells.Item(20,9).Value = ''
$ws.Cells.Item(20,10).Value = ''
$ws.Cells.Item(20,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(20,12).Value = 'Efsanevi metal kılıf güç transistörü.'
$ws.Cells.Item(20,13).Value = 'BASE,COLL,EMIT'
$ws.Cells.Item(20,14).Value = 20
$ws.Cells.Item(20,15).Value = ''

# Row 21: NE555
$ws.Cells.Item(21,1).Value = 'NE555'
$ws.Cells.Item(21,2).Value = 'IC'
$ws.Cells.Item(21,3).Value = 'Timer'
$ws.Cells.Item(21,4).Value = 'DIP-8'
$ws.Cells.Item(21,5).Value = '1-8'
$ws.Cells.Item(21,6).Value = 16
$ws.Cells.Item(21,7).Value = 0.2
$ws.Cells.Item(21,8).Value = 0.6
$ws.Cells.Item(21,9).Value = ''
$ws.Cells.Item(21,10).Value = ''
$ws.Cells.Item(21,11).Value = 'TEST_IC'
$ws.Cells.Item(21,12).Value = 'Hassas zamanlayıcı. Osilatör, PWM ve Timer devrelerinin kalbi.'
$ws.Cells.Item(21,13).Value = 'GND,TRIG,OUT,RST,CTRL,THR,DIS,VCC'
$ws.Cells.Item(21,14).Value = ''
$ws.Cells.Item(21,15).Value = ''

# Row 22: LM358
$ws.Cells.Item(22,1).Value = 'LM358'
$ws.Cells.Item(22,2).Value = 'IC'
$ws.Cells.Item(22,3).Value = 'OpAmp Dual'
$ws.Cells.Item(22,4).Value = 'DIP-8'
$ws.Cells.Item(22,5).Value = '1-8'
$ws.Cells.Item(22,6).Value = 32
$ws.Cells.Item(22,7).Value = 0.02
$ws.Cells.Item(22,8).Value = 0.5
$ws.Cells.Item(22,9).Value = ''
$ws.Cells.Item(22,10).Value = ''
$ws.Cells.Item(22,11).Value = 'TEST_IC'
$ws.Cells.Item(22,12).Value = 'Çift kanallı, tek beslemeyle çalışabilen genel amaçlı Op-Amp.'
$ws.Cells.Item(22,13).Value = 'OUT A,IN- A,IN+ A,GND,IN+ B,IN- B,OUT B,VCC'
$ws.Cells.Item(22,14).Value = ''
$ws.Cells.Item(22,15).Value = ''

# Row 23: LM741
$ws.Cells.Item(23,1).Value = 'LM741'
$ws.Cells.Item(23,2).Value = 'IC'
$ws.Cells.Item(23,3).Value = 'OpAmp Single'
$ws.Cells.Item(23,4).Value = 'DIP-8'
$ws.Cells.Item(23,5).Value = '1-8'
$ws.Cells.Item(23,6).Value = 22
$ws.Cells.Item(23,7).Value = 0.002
$ws.Cells.Item(23,8).Value = 0.5
$ws.Cells.Item(23,9).Value = ''
$ws.Cells.Item(23,10).Value = ''
$ws.Cells.Item(23,11).Value = 'TEST_IC'
$ws.Cells.Item(23,12).Value = 'Klasik tekli Op-Amp. Eğitim ve temel uygulamalar için.'
$ws.Cells.Item(23,13).Value = 'OFF,IN-,IN+,V-,OFF,OUT,V+,NC'
$ws.Cells.Item(23,14).Value = ''
$ws.Cells.Item(23,15).Value = ''

# Row 24: TL072
$ws.Cells.Item(24,1).Value = 'TL072'
$ws.Cells.Item(24,2).Value = 'IC'
$ws.Cells.Item(24,3).Value = 'JFET OpAmp'
$ws.Cells.Item(24,4).Value = 'DIP-8'
$ws.Cells.Item(24,5).Value = '1-8'
$ws.Cells.Item(24,6).Value = 36
$ws.Cells.Item(24,7).Value = 0.01
$ws.Cells.Item(24,8).Value = 0.6
$ws.Cells.Item(24,9).Value = ''
$ws.Cells.Item(24,10).Value = ''
$ws.Cells.Item(24,11).Value = 'TEST_IC'
$ws.Cells.Item(24,12).Value = 'Düşük gürültülü JFET girişli Op-Amp. Ses devreleri için ideal.'
$ws.Cells.Item(24,13).Value = 'OUT A,IN- A,IN+ A,V-,IN+ B,IN- B,OUT B,V+'
$ws.Cells.Item(24,14).Value = ''
$ws.Cells.Item(24,15).Value = ''

# Row 25: ULN2003
$ws.Cells.Item(25,1).Value = 'ULN2003'
$ws.Cells.Item(25,2).Value = 'IC'
$ws.Cells.Item(25,3).Value = 'Darlington'
$ws.Cells.Item(25,4).Value = 'DIP-16'
$ws.Cells.Item(25,5).Value = '1-16'
$ws.Cells.Item(25,6).Value = 50
$ws.Cells.Item(25,7).Value = 0.5
$ws.Cells.Item(25,8).Value = 1
$ws.Cells.Item(25,9).Value = ''
$ws.Cells.Item(25,10).Value = ''
$ws.Cells.Item(25,11).Value = 'TEST_IC'
$ws.Cells.Item(25,12).Value = '7 Kanal Darlington dizisi. Röle ve step motor sürmek için kullanılır.'
$ws.Cells.Item(25,13).Value = 'IN1,IN2,IN3,IN4,IN5,IN6,IN7,GND,COM,OUT7,OUT6,OUT5,OUT4,OUT3,OUT2,OUT1'
$ws.Cells.Item(25,14).Value = ''
$ws.Cells.Item(25,15).Value = ''

# Row 26: L293D
$ws.Cells.Item(26,1).Value = 'L293D'
$ws.Cells.Item(26,2).Value = 'IC'
$ws.Cells.Item(26,3).Value = 'Motor Driver'
$ws.Cells.Item(26,4).Value = 'DIP-16'
$ws.Cells.Item(26,5).Value = '1-16'
$ws.Cells.Item(26,6).Value = 36
$ws.Cells.Item(26,7).Value = 0.6
$ws.Cells.Item(26,8).Value = 1.5
$ws.Cells.Item(26,9).Value = ''
$ws.Cells.Item(26,10).Value = ''
$ws.Cells.Item(26,11).Value = 'TEST_IC'
$ws.Cells.Item(26,12).Value = 'Çift H-Köprüsü Motor Sürücü. DC motorları ileri-geri sürebilir.'
$ws.Cells.Item(26,13).Value = 'EN1,IN1,OUT1,GND,GND,OUT2,IN2,VCC2,EN2,IN3,OUT3,GND,GND,OUT4,IN4,VCC1'
$ws.Cells.Item(26,14).Value = ''
$ws.Cells.Item(26,15).Value = ''

# Row 27: CD4017
$ws.Cells.Item(27,1).Value = 'CD4017'
$ws.Cells.Item(27,2).Value = 'IC'
$ws.Cells.Item(27,3).Value = 'Counter'
$ws.Cells.Item(27,4).Value = 'DIP-16'
$ws.Cells.Item(27,5).Value = '1-16'
$ws.Cells.Item(27,6).Value = 15
$ws.Cells.Item(27,7).Value = 0.01
$ws.Cells.Item(27,8).Value = 0.5
$ws.Cells.Item(27,9).Value = ''
$ws.Cells.Item(27,10).Value = ''
$ws.Cells.Item(27,11).Value = 'TEST_IC'
$ws.Cells.Item(27,12).Value = 'Onlu sayıcı (Decade Counter). Yürüyen ışık devrelerinde popülerdir.'
$ws.Cells.Item(27,13).Value = '5,1,0,2,6,7,3,GND,8,4,9,CARRY,EN,CLK,RST,VCC'
$ws.Cells.Item(27,14).Value = ''
$ws.Cells.Item(27,15).Value = ''

# Row 28: PC817
$ws.Cells.Item(28,1).Value = 'PC817'
$ws.Cells.Item(28,2).Value = 'IC'
$ws.Cells.Item(28,3).Value = 'Optocoupler'
$ws.Cells.Item(28,4).Value = 'DIP-4'
$ws.Cells.Item(28,5).Value = 'AKEC'
$ws.Cells.Item(28,6).Value = 35
$ws.Cells.Item(28,7).Value = 0.05
$ws.Cells.Item(28,8).Value = 0.15
$ws.Cells.Item(28,9).Value = ''
$ws.Cells.Item(28,10).Value = ''
$ws.Cells.Item(28,11).Value = 'TEST_DIODE'
$ws.Cells.Item(28,12).Value = '4 Pinli Optokuplör. Sinyal izolasyonu sağlar.'
$ws.Cells.Item(28,13).Value = 'ANODE,CATHODE,EMITTER,COLLECTOR'
$ws.Cells.Item(28,14).Value = ''
$ws.Cells.Item(28,15).Value = ''

# Row 29: L7805
$ws.Cells.Item(29,1).Value = 'L7805'
$ws.Cells.Item(29,2).Value = 'REGULATOR'
$ws.Cells.Item(29,3).Value = 'Linear'
$ws.Cells.Item(29,4).Value = 'TO-220'
$ws.Cells.Item(29,5).Value = 'IGO'
$ws.Cells.Item(29,6).Value = 35
$ws.Cells.Item(29,7).Value = 1.5
$ws.Cells.Item(29,8).Value = 15
$ws.Cells.Item(29,9).Value = 5
$ws.Cells.Item(29,10).Value = ''
$ws.Cells.Item(29,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(29,12).Value = 'Pozitif 5V Sabit Regülatör.'
$ws.Cells.Item(29,13).Value = 'INPUT,GND,OUTPUT'
$ws.Cells.Item(29,14).Value = ''
$ws.Cells.Item(29,15).Value = ''

# Row 30: L7809
$ws.Cells.Item(30,1).Value = 'L7809'
$ws.Cells.Item(30,2).Value = 'REGULATOR'
$ws.Cells.Item(30,3).Value = 'Linear'
$ws.Cells.Item(30,4).Value = 'TO-220'
$ws.Cells.Item(30,5).Value = 'IGO'
$ws.Cells.Item(30,6).Value = 35
$ws.Cells.Item(30,7).Value = 1.5
$ws.Cells.Item(30,8).Value = 15
$ws.Cells.Item(30,9).Value = 9
$ws.Cells.Item(30,10).Value = ''
$ws.Cells.Item(30,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(30,12).Value = 'Pozitif 9V Sabit Regülatör.'
$ws.Cells.Item(30,13).Value = 'INPUT,GND,OUTPUT'
$ws.Cells.Item(30,14).Value = ''
$ws.Cells.Item(30,15).Value = ''

# Row 31: L7812
$ws.Cells.Item(31,1).Value = 'L7812'
$ws.Cells.Item(31,2).Value = 'REGULATOR'
$ws.Cells.Item(31,3).Value = 'Linear'
$ws.Cells.Item(31,4).Value = 'TO-220'
$ws.Cells.Item(31,5).Value = 'IGO'
$ws.Cells.Item(31,6).Value = 35
$ws.Cells.Item(31,7).Value = 1.5
$ws.Cells.Item(31,8).Value = 15
$ws.Cells.Item(31,9).Value = 12
$ws.Cells.Item(31,10).Value = ''
$ws.Cells.Item(31,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(31,12).Value = 'Pozitif 12V Sabit Regülatör.'
$ws.Cells.Item(31,13).Value = 'INPUT,GND,OUTPUT'
$ws.Cells.Item(31,14).Value = ''
$ws.Cells.Item(31,15).Value = ''

# Row 32: L7905
$ws.Cells.Item(32,1).Value = 'L7905'
$ws.Cells.Item(32,2).Value = 'REGULATOR'
$ws.Cells.Item(32,3).Value = 'Negative'
$ws.Cells.Item(32,4).Value = 'TO-220'
$ws.Cells.Item(32,5).Value = 'GIO'
$ws.Cells.Item(32,6).Value = -35
$ws.Cells.Item(32,7).Value = 1.5
$ws.Cells.Item(32,8).Value = 15
$ws.Cells.Item(32,9).Value = -5
$ws.Cells.Item(32,10).Value = ''
$ws.Cells.Item(32,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(32,12).Value = 'Negatif -5V Sabit Regülatör. (GND-Input-Out).'
$ws.Cells.Item(32,13).Value = 'GND,INPUT,OUTPUT'
$ws.Cells.Item(32,14).Value = ''
$ws.Cells.Item(32,15).Value = ''

# Row 33: LM317
$ws.Cells.Item(33,1).Value = 'LM317'
$ws.Cells.Item(33,2).Value = 'REGULATOR'
$ws.Cells.Item(33,3).Value = 'Adjust'
$ws.Cells.Item(33,4).Value = 'TO-220'
$ws.Cells.Item(33,5).Value = 'AOI'
$ws.Cells.Item(33,6).Value = 40
$ws.Cells.Item(33,7).Value = 1.5
$ws.Cells.Item(33,8).Value = 20
$ws.Cells.Item(33,9).Value = 1.25
$ws.Cells.Item(33,10).Value = ''
$ws.Cells.Item(33,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(33,12).Value = 'Ayarlanabilir Pozitif Regülatör (1.2V - 37V).'
$ws.Cells.Item(33,13).Value = 'ADJ,OUTPUT,INPUT'
$ws.Cells.Item(33,14).Value = ''
$ws.Cells.Item(33,15).Value = ''

# Row 34: AMS1117-3.3
$ws.Cells.Item(34,1).Value = 'AMS1117-3.3'
$ws.Cells.Item(34,2).Value = 'REGULATOR'
$ws.Cells.Item(34,3).Value = 'LDO'
$ws.Cells.Item(34,4).Value = 'SOT-223'
$ws.Cells.Item(34,5).Value = 'GOI'
$ws.Cells.Item(34,6).Value = 15
$ws.Cells.Item(34,7).Value = 0.8
$ws.Cells.Item(34,8).Value = 1
$ws.Cells.Item(34,9).Value = 3.3
$ws.Cells.Item(34,10).Value = ''
$ws.Cells.Item(34,11).Value = 'TEST_REGULATOR'
$ws.Cells.Item(34,12).Value = '3.3V LDO Regülatör (SMD).'
$ws.Cells.Item(34,13).Value = 'GND,OUTPUT,INPUT'
$ws.Cells.Item(34,14).Value = ''
$ws.Cells.Item(34,15).Value = ''

# Row 35: 1N4007
$ws.Cells.Item(35,1).Value = '1N4007'
$ws.Cells.Item(35,2).Value = 'DIODE'
$ws.Cells.Item(35,3).Value = 'Standard'
$ws.Cells.Item(35,4).Value = 'DO-41'
$ws.Cells.Item(35,5).Value = 'AK'
$ws.Cells.Item(35,6).Value = 1000
$ws.Cells.Item(35,7).Value = 1
$ws.Cells.Item(35,8).Value = 0
$ws.Cells.Item(35,9).Value = ''
$ws.Cells.Item(35,10).Value = ''
$ws.Cells.Item(35,11).Value = 'TEST_DIODE'
$ws.Cells.Item(35,12).Value = 'Genel amaçlı doğrultucu diyot.'
$ws.Cells.Item(35,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(35,14).Value = ''
$ws.Cells.Item(35,15).Value = ''

# Row 36: 1N4148
$ws.Cells.Item(36,1).Value = '1N4148'
$ws.Cells.Item(36,2).Value = 'DIODE'
$ws.Cells.Item(36,3).Value = 'Switching'
$ws.Cells.Item(36,4).Value = 'DO-35'
$ws.Cells.Item(36,5).Value = 'AK'
$ws.Cells.Item(36,6).Value = 100
$ws.Cells.Item(36,7).Value = 0.2
$ws.Cells.Item(36,8).Value = 0
$ws.Cells.Item(36,9).Value = ''
$ws.Cells.Item(36,10).Value = ''
$ws.Cells.Item(36,11).Value = 'TEST_DIODE'
$ws.Cells.Item(36,12).Value = 'Yüksek hızlı sinyal diyodu.'
$ws.Cells.Item(36,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(36,14).Value = ''
$ws.Cells.Item(36,15).Value = 4

# Row 37: UF4007
$ws.Cells.Item(37,1).Value = 'UF4007'
$ws.Cells.Item(37,2).Value = 'DIODE'
$ws.Cells.Item(37,3).Value = 'Fast Rec'
$ws.Cells.Item(37,4).Value = 'DO-41'
$ws.Cells.Item(37,5).Value = 'AK'
$ws.Cells.Item(37,6).Value = 1000
$ws.Cells.Item(37,7).Value = 1
$ws.Cells.Item(37,8).Value = 0
$ws.Cells.Item(37,9).Value = ''
$ws.Cells.Item(37,10).Value = ''
$ws.Cells.Item(37,11).Value = 'TEST_DIODE'
$ws.Cells.Item(37,12).Value = 'Ultra Hızlı (Fast Recovery) diyot.'
$ws.Cells.Item(37,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(37,14).Value = ''
$ws.Cells.Item(37,15).Value = 75

# Row 38: 1N5819
$ws.Cells.Item(38,1).Value = '1N5819'
$ws.Cells.Item(38,2).Value = 'DIODE'
$ws.Cells.Item(38,3).Value = 'Schottky'
$ws.Cells.Item(38,4).Value = 'DO-41'
$ws.Cells.Item(38,5).Value = 'AK'
$ws.Cells.Item(38,6).Value = 40
$ws.Cells.Item(38,7).Value = 1
$ws.Cells.Item(38,8).Value = 0
$ws.Cells.Item(38,9).Value = ''
$ws.Cells.Item(38,10).Value = ''
$ws.Cells.Item(38,11).Value = 'TEST_DIODE'
$ws.Cells.Item(38,12).Value = 'Düşük voltaj düşümlü Schottky diyot.'
$ws.Cells.Item(38,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(38,14).Value = ''
$ws.Cells.Item(38,15).Value = ''

# Row 39: BC846
$ws.Cells.Item(39,1).Value = 'BC846'
$ws.Cells.Item(39,2).Value = 'BJT'
$ws.Cells.Item(39,3).Value = 'NPN'
$ws.Cells.Item(39,4).Value = 'SOT-23'
$ws.Cells.Item(39,5).Value = 'CBE'
$ws.Cells.Item(39,6).Value = 65
$ws.Cells.Item(39,7).Value = 0.1
$ws.Cells.Item(39,8).Value = 0.25
$ws.Cells.Item(39,9).Value = ''
$ws.Cells.Item(39,10).Value = ''
$ws.Cells.Item(39,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(39,12).Value = 'SMD NPN Transistör.'
$ws.Cells.Item(39,13).Value = 'BASE,EMIT,COLL'
$ws.Cells.Item(39,14).Value = 110
$ws.Cells.Item(39,15).Value = ''

# Row 40: MMBT3904
$ws.Cells.Item(40,1).Value = 'MMBT3904'
$ws.Cells.Item(40,2).Value = 'BJT'
$ws.Cells.Item(40,3).Value = 'NPN'
$ws.Cells.Item(40,4).Value = 'SOT-23'
$ws.Cells.Item(40,5).Value = 'EBC'
$ws.Cells.Item(40,6).Value = 40
$ws.Cells.Item(40,7).Value = 0.2
$ws.Cells.Item(40,8).Value = 0.35
$ws.Cells.Item(40,9).Value = ''
$ws.Cells.Item(40,10).Value = ''
$ws.Cells.Item(40,11).Value = 'TEST_BJT_NPN'
$ws.Cells.Item(40,12).Value = '2N3904''ün SMD versiyonu.'
$ws.Cells.Item(40,13).Value = 'BASE,EMIT,COLL'
$ws.Cells.Item(40,14).Value = 100
$ws.Cells.Item(40,15).Value = ''

# Row 41: M7
$ws.Cells.Item(41,1).Value = 'M7'
$ws.Cells.Item(41,2).Value = 'DIODE'
$ws.Cells.Item(41,3).Value = 'Standard'
$ws.Cells.Item(41,4).Value = 'SMA'
$ws.Cells.Item(41,5).Value = 'AK'
$ws.Cells.Item(41,6).Value = 1000
$ws.Cells.Item(41,7).Value = 1
$ws.Cells.Item(41,8).Value = 0
$ws.Cells.Item(41,9).Value = ''
$ws.Cells.Item(41,10).Value = ''
$ws.Cells.Item(41,11).Value = 'TEST_DIODE'
$ws.Cells.Item(41,12).Value = '1N4007''nin SMD versiyonu.'
$ws.Cells.Item(41,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(41,14).Value = ''
$ws.Cells.Item(41,15).Value = ''

# Row 42: SS14
$ws.Cells.Item(42,1).Value = 'SS14'
$ws.Cells.Item(42,2).Value = 'DIODE'
$ws.Cells.Item(42,3).Value = 'Schottky'
$ws.Cells.Item(42,4).Value = 'SMA'
$ws.Cells.Item(42,5).Value = 'AK'
$ws.Cells.Item(42,6).Value = 40
$ws.Cells.Item(42,7).Value = 1
$ws.Cells.Item(42,8).Value = 0
$ws.Cells.Item(42,9).Value = ''
$ws.Cells.Item(42,10).Value = ''
$ws.Cells.Item(42,11).Value = 'TEST_DIODE'
$ws.Cells.Item(42,12).Value = '1N5819''un SMD versiyonu.'
$ws.Cells.Item(42,13).Value = 'ANODE,CATHODE'
$ws.Cells.Item(42,14).Value = ''
$ws.Cells.Item(42,15).Value = ''

# --- SMDCodes sheet updates ---
$ws2 = $wb.Worksheets.Item("SMDCodes")
$ws2.Cells.Item(6,1).Value = 'M7'
$ws2.Cells.Item(6,2).Value = 'M7'
$ws2.Cells.Item(6,3).Value = 'SMA'
$ws2.Cells.Item(7,1).Value = 'SS14'
$ws2.Cells.Item(7,2).Value = 'SS14'
$ws2.Cells.Item(7,3).Value = 'SMA'
$ws2.Cells.Item(8,1).Value = 'J3Y'
$ws2.Cells.Item(8,2).Value = 'S8050'
$ws2.Cells.Item(8,3).Value = 'SOT-23'

# --- Packages sheet updates ---
$ws3 = $wb.Worksheets.Item("Packages")
$ws3.Cells.Item(12,1).Value = 'DIP-4'
$ws3.Cells.Item(12,2).Value = 'assets/packages/dip-4.png'
$ws3.Cells.Item(13,1).Value = 'DIP-8'
$ws3.Cells.Item(13,2).Value = 'assets/packages/dip-8.png'
$ws3.Cells.Item(14,1).Value = 'DIP-14'
$ws3.Cells.Item(14,2).Value = 'assets/packages/dip-14.png'
$ws3.Cells.Item(15,1).Value = 'DIP-16'
$ws3.Cells.Item(15,2).Value = 'assets/packages/dip-16.png'
